# Fix 3.2 beta3: excel template now only exports the first image, and the
# barcode cell's canned greeting text is dropped in favor of the bare
# {BARCODE} placeholder. A new "图片" / {IMAGE} column is added to the
# report template so the generated report has a slot for the image.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 used to read "你好 {BARCODE}" - strip the greeting, keep the placeholder.
$ws.Range("A2").Value = "{BARCODE}"

# New column D: header "图片" (styled like the other header cells in row 1)
# and placeholder "{IMAGE}" (styled like the other value cells in row 2).
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "图片"

$ws.Range("C2").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("D2").Value = "{IMAGE}"

# Match the selection/active cell that Excel leaves after adding the column.
$null = $ws.Range("D3").Select()
